# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on "Weekly Quantity" to "Weekly_PO_Qty"
# 2. Rename the "Requested quantity" header on "Monthly Trend" to "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet (after the existing sheets) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename column B header ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename column B header ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header style (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$wsForecast.Cells.Item(2,1).Value = 45130.99999999999
$wsForecast.Cells.Item(2,2).Value = 1
$wsForecast.Cells.Item(2,3).Value = 0.7877744799399825
$wsForecast.Cells.Item(2,4).Value = 1.929812484749017
$wsForecast.Cells.Item(3,1).Value = 45144.99999999999
$wsForecast.Cells.Item(3,2).Value = 1
$wsForecast.Cells.Item(3,3).Value = 0.7666247660539964
$wsForecast.Cells.Item(3,4).Value = 1.900812693516746
$wsForecast.Cells.Item(4,1).Value = 45158.99999999999
$wsForecast.Cells.Item(4,2).Value = 1
$wsForecast.Cells.Item(4,3).Value = 0.7547043364503723
$wsForecast.Cells.Item(4,4).Value = 1.911739042112634
$wsForecast.Cells.Item(5,1).Value = 45165.99999999999
$wsForecast.Cells.Item(5,2).Value = 1
$wsForecast.Cells.Item(5,3).Value = 0.7471787685600931
$wsForecast.Cells.Item(5,4).Value = 1.874090846118089
$wsForecast.Cells.Item(6,1).Value = 45172.99999999999
$wsForecast.Cells.Item(6,2).Value = 1
$wsForecast.Cells.Item(6,3).Value = 0.7523958125707086
$wsForecast.Cells.Item(6,4).Value = 1.92077967076245
$wsForecast.Cells.Item(7,1).Value = 45179.99999999999
$wsForecast.Cells.Item(7,2).Value = 1
$wsForecast.Cells.Item(7,3).Value = 0.781442611041673
$wsForecast.Cells.Item(7,4).Value = 1.965392440278622
$wsForecast.Cells.Item(8,1).Value = 45186.99999999999
$wsForecast.Cells.Item(8,2).Value = 1
$wsForecast.Cells.Item(8,3).Value = 0.7412012079919518
$wsForecast.Cells.Item(8,4).Value = 1.886587120870039
$wsForecast.Cells.Item(9,1).Value = 45193.99999999999
$wsForecast.Cells.Item(9,2).Value = 1
$wsForecast.Cells.Item(9,3).Value = 0.7244602590250748
$wsForecast.Cells.Item(9,4).Value = 1.939230162129078
$wsForecast.Cells.Item(10,1).Value = 45200.99999999999
$wsForecast.Cells.Item(10,2).Value = 1
$wsForecast.Cells.Item(10,3).Value = 0.7378536922620049
$wsForecast.Cells.Item(10,4).Value = 1.917196815926371
$wsForecast.Cells.Item(11,1).Value = 45207.99999999999
$wsForecast.Cells.Item(11,2).Value = 1
$wsForecast.Cells.Item(11,3).Value = 0.704313358298723
$wsForecast.Cells.Item(11,4).Value = 1.880196363353363
$wsForecast.Cells.Item(12,1).Value = 45214.99999999999
$wsForecast.Cells.Item(12,2).Value = 1
$wsForecast.Cells.Item(12,3).Value = 0.7196842380849118
$wsForecast.Cells.Item(12,4).Value = 1.916752214113955

# Match the date style (column A) used on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)

# Put selection back on A1 and leave the first sheet active, matching source workbook state
$wsForecast.Range("A1").Select()
$wsWeekly.Activate()
